$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values could be misread as numbers by Excel;
# force Text format, assign, then restore default style so no style
# attributes are left behind (matches original unstyled cells).
$numericLookingCells = @(
    "D5", "D6", "D10", "D11", "D12", "D14", "D16", "D20", "D23", "D24", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D39", "D40", "D42", "D44", "D46", "D47", "D48", "D49"
)
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.731.91"
$ws.Range("E2").Value = "  +4.29%  "
$ws.Range("D3").Value = "2.264.71"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "303.78"
$ws.Range("E5").Value = "  +3.17%  "
$ws.Range("D6").Value = "91.47"
$ws.Range("E6").Value = "  +4.69%  "
$ws.Range("E7").Value = "  +3.35%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("D10").Value = "32.27"
$ws.Range("E10").Value = "  +5.12%  "
$ws.Range("D11").Value = "53.27"
$ws.Range("E11").Value = "  +3.75%  "
$ws.Range("D12").Value = "0.0795"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").Value = "6.58"
$ws.Range("E14").Value = "  +3.00%  "
$ws.Range("D15").Value = "2.617.47"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "14.18"
$ws.Range("E16").Value = "  +2.73%  "
$ws.Range("D17").Value = "2.267.95"
$ws.Range("E17").Value = "  +4.35%  "
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("D19").Value = "41.657.78"
$ws.Range("D20").Value = "12.39"
$ws.Range("E20").Value = "  +10.64%  "
$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("E22").Value = "  +2.59%  "
$ws.Range("D23").Value = "66.62"
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").Value = "240.98"
$ws.Range("E25").Value = "  +4.51%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +5.29%  "
$ws.Range("D28").Value = "24.11"
$ws.Range("E28").Value = "  +4.41%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "9.49"
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.07"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").Value = "34.71"
$ws.Range("E31").Value = "  +9.69%  "
$ws.Range("D32").Value = "160.51"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "5.15"
$ws.Range("E34").Value = "  +4.05%  "
$ws.Range("D35").Value = "0.0744"
$ws.Range("E35").Value = "  +4.51%  "
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("E37").Value = "  +2.05%  "
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("D39").Value = "16.61"
$ws.Range("E39").Value = "  +7.14%  "
$ws.Range("D40").Value = "0.104"
$ws.Range("E40").Value = "  +4.11%  "
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("D42").Value = "3.91"
$ws.Range("E42").Value = "  +4.83%  "
$ws.Range("D43").Value = "2.062.17"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").Value = "19.27"
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("E45").Value = "  +2.65%  "
$ws.Range("D46").Value = "10.13"
$ws.Range("E46").Value = "  +2.49%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "2.86"
$ws.Range("E47").Value = "  +3.75%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "2.04"
$ws.Range("E48").Value = "  +5.28%  "
$ws.Range("D49").Value = "72.82"
$ws.Range("E49").Value = "  +8.08%  "
$ws.Range("E50").Value = "  +4.17%  "
$ws.Range("E51").Value = "  +3.24%  "

foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
